# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Integral"     (bound to the slide master / visible slides)
#   ppt/theme/theme2.xml  -> "Office Theme" (bound to the notes master)
#
# The authored change swaps the two themes' contents: the slide master
# (theme1.xml) becomes the default "Office Theme" palette, and the notes
# master (theme2.xml) becomes the old "Integral" palette. The font scheme
# and format scheme are identical between the two themes already, so the
# only substantive difference is the 12-slot colour scheme.
#
# Apply the "Office Theme" colour scheme to the presentation's slide master
# (theme1.xml) via the standard ColorScheme COM surface. Colour indices
# follow the classic ppColorSchemeIndex ordering:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4
#   9 accent5  10 accent6  11 hlink  12 folHlink
# RGBColor.RGB uses the OLE COLORREF (0x00BBGGRR) byte order.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
